# "fixxed add bird window" - append the missing CageDB row (row 44) that was
# dropped from Sheet1: CageID 250, Length 22.55, Width 22, Height 22,
# Material "Plastic".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A44").Value = 250
$ws.Range("B44").Value = 22.55
$ws.Range("C44").Value = 22
$ws.Range("D44").Value = 22
$ws.Range("E44").Value = "Plastic"
